# Generate Report for Handback
# The file "4a88754c-d4d3-40bf-883e-05388a36cbfc.md" transitions from
# "Ready for handoff" to "Handed back: in sync with en-US" across all sheets,
# and the Latest Handback DateTime is refreshed on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the 4a88754c-...md file
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: row 3 is the 4a88754c-...md file
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("H3").Value = "2016-03-22 06:34:27"

# de-de sheet: row 3 is the 4a88754c-...md file
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("H3").Value = "2016-03-22 06:34:40"
